# Card2 ("Machine_Service_Lookup.xlsx"): add a new 'Correction ' column (N)
# to the right of the existing 'Event' column (M), matching the layout
# already used by the other Card sheets (e.g. Card24, which has
# ... | L:Date | M:Event | N:Correction | O:Servised by).
#
# Changes required:
#   - Dimension grows from A1:M13 to A1:N13
#   - M1 header text "Event " (trailing space) -> "Event" (no trailing space)
#   - N1 new header "Correction " (with trailing space), same style as M1
#   - M2:M13 (currently blank) get the literal text "nan"
#   - N2:N13 are created as new, still-blank cells under the new column

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card2")

# --- Header row -----------------------------------------------------
# Trim the trailing space from the existing "Event " header.
$ws.Range("M1").Value = "Event"

# Create N1 with the same formatting as M1 (bold/border/centered header
# style), then set its text to the new column name.
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("N1").Value = "Correction "

# --- Data rows --------------------------------------------------------
for ($r = 2; $r -le 13; $r++) {
    # M column: previously-blank cells now hold the literal text "nan".
    $ws.Cells.Item($r, 13).Value = "nan"

    # N column: brand-new column, cells stay empty - just "touch" the
    # cell (no-op format write) so it is materialized as part of the
    # sheet/used-range rather than left completely absent.
    $ws.Cells.Item($r, 14).Font.Bold = $false
}
